$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.667069666666666
$ws.Range("H2").Value = 11.001209
$ws.Range("I2").Value = 0.01298011522000835
$ws.Range("J2").Value = 0.01298011522000835
$ws.Range("M2").Value = 20.88869433333333
$ws.Range("N2").Value = 62.666083
$ws.Range("O2").Value = 0.06073223131780172
$ws.Range("P2").Value = 0.06073223131780172
$ws.Range("Q2").Value = 76.60029736603855
$ws.Range("R2").Value = 689.402676294347
$ws.Range("S2").Value = 0.000788311360073266
$ws.Range("T2").Value = 0.0007883113600732661
$ws.Range("G3").Value = 3.667069666666666
$ws.Range("H3").Value = 11.001209
$ws.Range("I3").Value = 0.01298011522000835
$ws.Range("J3").Value = 0.01298011522000835
$ws.Range("O3").Value = 0.1122209945765712
$ws.Range("P3").Value = 0.1122209945765712
$ws.Range("Q3").Value = 141.5420011541428
$ws.Range("R3").Value = 1273.878010387285
$ws.Range("S3").Value = 0.001456641439707827
$ws.Range("T3").Value = 0.001456641439707827
$ws.Range("G4").Value = 3.667069666666666
$ws.Range("H4").Value = 11.001209
$ws.Range("I4").Value = 0.01298011522000835
$ws.Range("J4").Value = 0.01298011522000835
$ws.Range("M4").Value = 277.3327996666667
$ws.Range("N4").Value = 831.9983990000001
$ws.Range("O4").Value = 0.8063232422570387
$ws.Range("P4").Value = 0.8063232422570388
$ws.Range("Q4").Value = 1016.998697229377
$ws.Range("R4").Value = 9152.988275064392
$ws.Range("S4").Value = 0.01046616858906707
$ws.Range("T4").Value = 0.01046616858906707
$ws.Range("G5").Value = 3.667069666666666
$ws.Range("H5").Value = 11.001209
$ws.Range("I5").Value = 0.01298011522000835
$ws.Range("J5").Value = 0.01298011522000835
$ws.Range("M5").Value = 7.127805333333332
$ws.Range("N5").Value = 21.383416
$ws.Range("O5").Value = 0.02072353184858837
$ws.Range("P5").Value = 0.02072353184858837
$ws.Range("Q5").Value = 26.13815872777155
$ws.Range("R5").Value = 235.243428549944
$ws.Range("S5").Value = 0.0002689938311601897
$ws.Range("T5").Value = 0.0002689938311601898
$ws.Range("I6").Value = 0.5954329572989919
$ws.Range("J6").Value = 0.595432957298992
$ws.Range("M6").Value = 20.88869433333333
$ws.Range("N6").Value = 62.666083
$ws.Range("O6").Value = 0.06073223131780172
$ws.Range("P6").Value = 0.06073223131780172
$ws.Range("Q6").Value = 3513.862613510234
$ws.Range("R6").Value = 31624.76352159211
$ws.Range("S6").Value = 0.03616197209692514
$ws.Range("T6").Value = 0.03616197209692514
$ws.Range("I7").Value = 0.5954329572989919
$ws.Range("J7").Value = 0.595432957298992
$ws.Range("O7").Value = 0.1122209945765712
$ws.Range("P7").Value = 0.1122209945765712
$ws.Range("S7").Value = 0.06682007867176194
$ws.Range("T7").Value = 0.06682007867176196
$ws.Range("I8").Value = 0.5954329572989919
$ws.Range("J8").Value = 0.595432957298992
$ws.Range("M8").Value = 277.3327996666667
$ws.Range("N8").Value = 831.9983990000001
$ws.Range("O8").Value = 0.8063232422570387
$ws.Range("P8").Value = 0.8063232422570388
$ws.Range("Q8").Value = 46652.47816345041
$ws.Range("R8").Value = 419872.3034710537
$ws.Range("S8").Value = 0.48011143267602
$ws.Range("T8").Value = 0.4801114326760202
$ws.Range("I9").Value = 0.5954329572989919
$ws.Range("J9").Value = 0.595432957298992
$ws.Range("M9").Value = 7.127805333333332
$ws.Range("N9").Value = 21.383416
$ws.Range("O9").Value = 0.02072353184858837
$ws.Range("P9").Value = 0.02072353184858837
$ws.Range("Q9").Value = 1199.027965917968
$ws.Range("R9").Value = 10791.25169326171
$ws.Range("S9").Value = 0.01233947385428482
$ws.Range("T9").Value = 0.01233947385428482
$ws.Range("G10").Value = 110.4727123333333
$ws.Range("H10").Value = 331.418137
$ws.Range("I10").Value = 0.3910338949346852
$ws.Range("J10").Value = 0.3910338949346853
$ws.Range("M10").Value = 20.88869433333333
$ws.Range("N10").Value = 62.666083
$ws.Range("O10").Value = 0.06073223131780172
$ws.Range("P10").Value = 0.06073223131780172
$ws.Range("Q10").Value = 2307.630720105264
$ws.Range("R10").Value = 20768.67648094737
$ws.Range("S10").Value = 0.02374836096027428
$ws.Range("T10").Value = 0.02374836096027428
$ws.Range("G11").Value = 110.4727123333333
$ws.Range("H11").Value = 331.418137
$ws.Range("I11").Value = 0.3910338949346852
$ws.Range("J11").Value = 0.3910338949346853
$ws.Range("O11").Value = 0.1122209945765712
$ws.Range("P11").Value = 0.1122209945765712
$ws.Range("Q11").Value = 4264.039191488667
$ws.Range("R11").Value = 38376.35272339801
$ws.Range("S11").Value = 0.04388221260272083
$ws.Range("T11").Value = 0.04388221260272084
$ws.Range("G12").Value = 110.4727123333333
$ws.Range("H12").Value = 331.418137
$ws.Range("I12").Value = 0.3910338949346852
$ws.Range("J12").Value = 0.3910338949346853
$ws.Range("M12").Value = 277.3327996666667
$ws.Range("N12").Value = 831.9983990000001
$ws.Range("O12").Value = 0.8063232422570387
$ws.Range("P12").Value = 0.8063232422570388
$ws.Range("Q12").Value = 30637.70659817363
$ws.Range("R12").Value = 275739.3593835627
$ws.Range("S12").Value = 0.3152997179961336
$ws.Range("T12").Value = 0.3152997179961337
$ws.Range("G13").Value = 110.4727123333333
$ws.Range("H13").Value = 331.418137
$ws.Range("I13").Value = 0.3910338949346852
$ws.Range("J13").Value = 0.3910338949346853
$ws.Range("M13").Value = 7.127805333333332
$ws.Range("N13").Value = 21.383416
$ws.Range("O13").Value = 0.02072353184858837
$ws.Range("P13").Value = 0.02072353184858837
$ws.Range("Q13").Value = 787.4279881573323
$ws.Range("R13").Value = 7086.851893415991
$ws.Range("S13").Value = 0.008103603375556505
$ws.Range("T13").Value = 0.008103603375556509
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1562396666666667
$ws.Range("H14").Value = 0.468719
$ws.Range("I14").Value = 0.0005530325463144183
$ws.Range("J14").Value = 0.0005530325463144184
$ws.Range("M14").Value = 20.88869433333333
$ws.Range("N14").Value = 62.666083
$ws.Range("O14").Value = 0.06073223131780172
$ws.Range("P14").Value = 0.06073223131780172
$ws.Range("Q14").Value = 3.263642639741889
$ws.Range("R14").Value = 29.372783757677
$ws.Range("S14").Value = 0.00003358690052904015
$ws.Range("T14").Value = 0.00003358690052904015
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1562396666666667
$ws.Range("H15").Value = 0.468719
$ws.Range("I15").Value = 0.0005530325463144183
$ws.Range("J15").Value = 0.0005530325463144184
$ws.Range("O15").Value = 0.1122209945765712
$ws.Range("P15").Value = 0.1122209945765712
$ws.Range("Q15").Value = 6.030557663159444
$ws.Range("R15").Value = 54.275018968435
$ws.Range("S15").Value = 0.00006206186238061771
$ws.Range("T15").Value = 0.00006206186238061772
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1562396666666667
$ws.Range("H16").Value = 0.468719
$ws.Range("I16").Value = 0.0005530325463144183
$ws.Range("J16").Value = 0.0005530325463144184
$ws.Range("M16").Value = 277.3327996666667
$ws.Range("N16").Value = 831.9983990000001
$ws.Range("O16").Value = 0.8063232422570387
$ws.Range("P16").Value = 0.8063232422570388
$ws.Range("Q16").Value = 43.33038417565345
$ws.Range("R16").Value = 389.973457580881
$ws.Range("S16").Value = 0.0004459229958179076
$ws.Range("T16").Value = 0.0004459229958179078
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1562396666666667
$ws.Range("H17").Value = 0.468719
$ws.Range("I17").Value = 0.0005530325463144183
$ws.Range("J17").Value = 0.0005530325463144184
$ws.Range("M17").Value = 7.127805333333332
$ws.Range("N17").Value = 21.383416
$ws.Range("O17").Value = 0.02072353184858837
$ws.Range("P17").Value = 0.02072353184858837
$ws.Range("Q17").Value = 1.113645929344889
$ws.Range("R17").Value = 10.022813364104
$ws.Range("S17").Value = 0.00001146078758685277
$ws.Range("T17").Value = 0.00001146078758685277